# Restore C10 and B11 on the "Rules" sheet to the values from the
# referenced revision:
#   - C10 goes from 18 back to 1
#   - B11 goes from the stray "1" label back to "R40"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
$ws.Range("B11").Value = "R40"
